$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2457.8333
$ws.Range("J17").Value = 2599.4546
$ws.Range("L17").Value = 7798.3638
$ws.Range("N17").Value = -8134.3638

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null

$ws.Range("H116").Value = 4486
$ws.Range("I116").Value = 4460
$ws.Range("K116").Value = 4460
$ws.Range("M116").Value = -1018

$ws.Range("H132").Value = 19703
$ws.Range("I132").Value = 19703
$ws.Range("K132").Value = 59109
$ws.Range("M132").Value = -56579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8462.6
$ws.Range("I31").Value = 8462.6
$ws.Range("K31").Value = 8462.6
$ws.Range("M31").Value = -8168.6

$ws.Range("H74").Value = 7451.2
$ws.Range("I74").Value = 7231.55
$ws.Range("J74").Value = 8329.799999999999
$ws.Range("K74").Value = 7231.55
$ws.Range("L74").Value = 8329.799999999999
$ws.Range("M74").Value = -6357.55
$ws.Range("N74").Value = -10077.8

$ws.Range("H77").Value = 7451.2
$ws.Range("I77").Value = 7231.55
$ws.Range("J77").Value = 8329.799999999999
$ws.Range("K77").Value = 36157.75
$ws.Range("L77").Value = 41649
$ws.Range("M77").Value = -31789.75
$ws.Range("N77").Value = -50385

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 132.52942
$ws.Range("I80").Value = 72.42856999999999
$ws.Range("J80").Value = 174.6
$ws.Range("K80").Value = 72.42856999999999
$ws.Range("L80").Value = 174.6
$ws.Range("M80").Value = 925.57143
$ws.Range("N80").Value = -2170.6

$ws.Range("H83").Value = 132.52942
$ws.Range("I83").Value = 72.42856999999999
$ws.Range("J83").Value = 174.6
$ws.Range("K83").Value = 362.14285
$ws.Range("L83").Value = 873
$ws.Range("M83").Value = 4629.85715
$ws.Range("N83").Value = -10857

$ws.Range("H94").Value = 666.5
$ws.Range("I94").Value = 633.3333
$ws.Range("K94").Value = 633.3333
$ws.Range("M94").Value = -182.3333

$ws.Range("H105").Value = 1427.5
$ws.Range("I105").Value = 1505
$ws.Range("K105").Value = 1505
$ws.Range("M105").Value = 242

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 742.7778
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 814.1667
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 814.1667
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -1040.1667

$ws.Range("H10").Value = 1850
$ws.Range("I10").Value = 1850
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1850
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1711
$ws.Range("N10").Value = $null

$ws.Range("H16").Value = 1666
$ws.Range("I16").Value = 1666
$ws.Range("J16").Value = 1666
$ws.Range("K16").Value = 1666
$ws.Range("L16").Value = 1666
$ws.Range("M16").Value = -1379
$ws.Range("N16").Value = -2240

$ws.Range("H33").Value = 743.44446
$ws.Range("I33").Value = 743.44446
$ws.Range("K33").Value = 743.44446
$ws.Range("M33").Value = -364.44446

$ws.Range("H105").Value = 1701.25
$ws.Range("I105").Value = 4015
$ws.Range("J105").Value = 930
$ws.Range("K105").Value = 4015
$ws.Range("L105").Value = 930
$ws.Range("M105").Value = -2268
$ws.Range("N105").Value = -4424

$ws.Range("H113").Value = 1666
$ws.Range("I113").Value = 1666
$ws.Range("J113").Value = 1666
$ws.Range("K113").Value = 1666
$ws.Range("L113").Value = 1666
$ws.Range("M113").Value = 504
$ws.Range("N113").Value = -6006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 195.25
$ws.Range("I23").Value = 123.5
$ws.Range("J23").Value = 231.125
$ws.Range("K23").Value = 370.5
$ws.Range("L23").Value = 693.375
$ws.Range("M23").Value = -135.5
$ws.Range("N23").Value = -1163.375

$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = $null

$ws.Range("H92").Value = 9500
$ws.Range("I92").Value = 9500
$ws.Range("K92").Value = 28500
$ws.Range("M92").Value = -27252

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 120.85714
$ws.Range("I13").Value = 20
$ws.Range("J13").Value = 128.61539
$ws.Range("K13").Value = 20
$ws.Range("L13").Value = 128.61539
$ws.Range("M13").Value = 119
$ws.Range("N13").Value = -406.61539

$ws.Range("H36").Value = 4998
$ws.Range("J36").Value = 4000
$ws.Range("L36").Value = 4000
$ws.Range("N36").Value = -4970

$ws.Range("H102").Value = 2036.1
$ws.Range("I102").Value = 2036.1
$ws.Range("K102").Value = 2036.1
$ws.Range("M102").Value = -414.0999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1348.5
$ws.Range("I10").Value = 4999
$ws.Range("K10").Value = 4999
$ws.Range("M10").Value = -4859

$ws.Range("H16").Value = 161.6
$ws.Range("I16").Value = 161.6
$ws.Range("K16").Value = 161.6
$ws.Range("M16").Value = 8.400000000000006

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null

$ws.Range("H64").Value = 57494.5
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null

$ws.Range("H67").Value = 57494.5
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

$ws.Range("H69").Value = 17708
$ws.Range("J69").Value = 21771
$ws.Range("L69").Value = 21771
$ws.Range("N69").Value = -23269

$ws.Range("H72").Value = 17708
$ws.Range("J72").Value = 21771
$ws.Range("L72").Value = 65313
$ws.Range("N72").Value = -72801

$ws.Range("H122").Value = 1068.1
$ws.Range("I122").Value = 1025.6666
$ws.Range("K122").Value = 3076.9998
$ws.Range("M122").Value = -626.9998000000001

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null
